{"js": "// The document has a single inline picture (a vehicular-access photo)\n// sitting in its own paragraph. The edit replaces that picture with a\n// text hyperlink whose visible text is the image's own URL, i.e. the\n// <w:drawing>/<wp:inline> run becomes a <w:hyperlink> run styled with\n// the \"Hyperlink\" character style.\n\nconst url =\n  \"https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Landed-Housing/TH03_Vehicular_Access.jpg?h=100%25&w=100%25\";\n\nconst pictures = context.document.body.inlinePictures;\npictures.load(\"items\");\nawait context.sync();\n\nif (pictures.items.length === 0) {\n  throw new Error(\"Expected an inline picture to convert into a hyperlink.\");\n}\n\n// There is exactly one inline picture in this document (the vehicular\n// access photo); replace it (and only it) with the hyperlink text.\nconst picture = pictures.items[0];\n\n// Replacing the picture's own range with text removes the drawing and\n// leaves the new text in its place, inside the same paragraph.\nconst linkRange = picture.insertText(url, Word.InsertLocation.replace);\n\n// Turning the inserted range into a hyperlink wraps it in <w:hyperlink>\n// and applies the built-in \"Hyperlink\" character style, matching the\n// target markup.\nlinkRange.hyperlink = url;\n\nawait context.sync();\n", "ps1": "# The document has a single inline picture (a vehicular-access photo)\n# living alone in its own paragraph. The edit swaps that picture for a\n# text hyperlink whose display text is the image's own URL, i.e. the\n# inline <w:drawing> run becomes a <w:hyperlink> run styled with the\n# built-in \"Hyperlink\" character style.\n\n$d = $word.ActiveDocument\n$url = \"https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Landed-Housing/TH03_Vehicular_Access.jpg?h=100%25&w=100%25\"\n\n# Walk backwards in case there is more than one picture; in this\n# document there is exactly one, the vehicular-access photo.\nfor ($i = $d.InlineShapes.Count; $i -ge 1; $i--) {\n    $shape = $d.InlineShapes.Item($i)\n    $rng = $shape.Range\n\n    # Deleting the shape removes the drawing but keeps the (now empty)\n    # range positioned where the picture used to be.\n    $shape.Delete()\n\n    # Put the URL text in its place, then turn that same range into a\n    # hyperlink pointing at the URL.\n    $rng.Text = $url\n    $d.Hyperlinks.Add($rng, $url)\n}\n"}
